$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "icon/config resource" column (B) that pairs with each
# effect-id row in column A. Row 1 is the header; row 2 keeps its own
# icon value, the remaining rows share the generic "Ssetting" resource.
$ws.Range("B1").Value = "Atlas_ResID"
$ws.Range("B2").Value = "msg_icon"
$ws.Range("B3:B15").Value = "Ssetting"

# Move the active selection, matching the saved-state cursor position.
[void]$ws.Range("E14").Select()
